$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# --- Update time_taken (column F) timestamps on the "data" sheet ---
# A re-run of the panel query re-stamped every gene row with a later time.
$newTimes = @(
    "2021-10-05 14:21:25.998756",
    "2021-10-05 14:21:25.998764",
    "2021-10-05 14:21:25.998767",
    "2021-10-05 14:21:25.998769",
    "2021-10-05 14:21:25.998772",
    "2021-10-05 14:21:25.998775",
    "2021-10-05 14:21:25.998777",
    "2021-10-05 14:21:25.998780",
    "2021-10-05 14:21:25.998783",
    "2021-10-05 14:21:25.998785"
)

for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $newTimes[$i]
}

# --- Add a new "metadata" worksheet directly after the "data" sheet ---
$meta = $wb.Worksheets.Add($null, $ws)
$meta.Name = "metadata"

# Header row (B1:G1)
$meta.Cells.Item(1, 2).Value = "data_name"
$meta.Cells.Item(1, 3).Value = "data_id"
$meta.Cells.Item(1, 4).Value = "data_version"
$meta.Cells.Item(1, 5).Value = "data_version_created"
$meta.Cells.Item(1, 6).Value = "panel_query_time"
$meta.Cells.Item(1, 7).Value = "panel_get_request"

# Match the bordered/bold "header" look used on the data sheet's header row
$ws.Range("B1:F1").Copy()
$meta.Range("B1:F1").PasteSpecial(-4122)
$ws.Range("B1").Copy()
$meta.Range("G1").PasteSpecial(-4122)

# Data row (row 2)
$meta.Cells.Item(2, 1).Value = 0
$meta.Cells.Item(2, 2).Value = "Lipoprotein lipase deficiency"
$meta.Cells.Item(2, 3).Value = 527
$meta.Cells.Item(2, 5).Value = "2021-08-05T15:48:31.031887Z"
$meta.Cells.Item(2, 6).Value = "2021-10-05 14:21:25.995502"
$meta.Cells.Item(2, 7).Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/527/?format=json"

# data_version must stay textual ("1.19"), not be coerced into a float
$meta.Cells.Item(2, 4).NumberFormat = "@"
$meta.Cells.Item(2, 4).Value = "1.19"
$meta.Cells.Item(2, 4).Style = "Normal"

# Copy the bordered/centered numeric-index look from the data sheet's A column
$ws.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)

$excel.CutCopyMode = 0
